$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 120, shifting existing rows 120:136 down to 121:137.
$ws.Rows.Item(120).Insert(-4121)

# Populate the newly inserted row 120 with the new price record.
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 45223
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112044
$ws.Cells.Item(120, 7).Value = "Perejil"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 150
$ws.Cells.Item(120, 11).Value = 1500
$ws.Cells.Item(120, 12).Value = 1500
$ws.Cells.Item(120, 13).Value = 1500
$ws.Cells.Item(120, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(120, 15).Value = "Región de Ñuble"
$ws.Cells.Item(120, 16).Value = 1500
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
